$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Your algorithm best solution obj" column (I) with the final
# benchmark results used for the release.
$ws.Range("I4").Value = 180.97476
$ws.Range("I5").Value = 46.11374
$ws.Range("I6").Value = 47.678894
$ws.Range("I7").Value = 10.917183
$ws.Range("I8").Value = 18.962624000000002
$ws.Range("I9").Value = 4.4145764999999999
$ws.Range("I10").Value = 14.834254

# Page setup used for the final release printout.
$ws.PageSetup.PaperSize = 9

# Leave the selection where the author last left it before saving.
$ws.Activate()
$null = $ws.Range("L6").Select()
